$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = $true
$ws.Range("B1").Value = "who discovered mongo park?"
$ws.Range("C1").Value = "a"
$ws.Range("D1").Value = "a: Julius beger;b: Akpan;c: James;d: Titus;"
$ws.Range("E1").Value = 9
$ws.Range("F1").Value = "natural sciences:mathematical sciences;"

$ws.Range("A2").Value = $true
$ws.Range("B2").Value = "who discovered mongo park?"
$ws.Range("C2").Value = "a"
$ws.Range("D2").Value = "a: Julius beger;b: Akpan;c: James;d: Titus;"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "natural sciences:mathematical sciences;"

$ws.Range("G5").Select()
